$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E16:E22) gets reversed: previously 2101..2107 top-to-bottom,
# now 2107..2101 top-to-bottom.
$periods = @("2107", "2106", "2105", "2104", "2103", "2102", "2101")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

# "Valor Mora" column: the values tied to the first and last rows are swapped.
$ws.Range("F16").Value = 58533
$ws.Range("F22").Value = 70240
